# #5: property aircraft done
# The "建物" (building) sheet's property_category column (I) was incorrectly
# tagged as "land" for every data row; relabel it as "building".
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("建物")

$lastRow = $ws.Cells.Item($ws.Rows.Count, 9).End(-4162).Row  # xlUp
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 9).Value = "building"
}
